$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cell (A2, merged A2:F2): drop the trailing "Date :13/09/2022"
# portion of the third rich-text run, keeping the "Store code : " and the
# bold "Mcd 006 - Fort" runs (and their formatting) intact.
$title = $ws.Range("A2")

$prefixLen = "Store code : Mcd 006 - Fort".Length   # 27
$tailChars = $title.Characters($prefixLen + 1, 1000)
$tailChars.Text = "               "

# Re-apply the original (non-bold, 17pt, Calibri) look to the run that
# replaces the old date text, since rewriting the characters resets run
# formatting.
$newTail = $title.Characters($prefixLen + 1, 15)
$newTail.Font.Name = "Calibri"
$newTail.Font.Size = 17
$newTail.Font.Bold = $false

# Make sure the bold store-name run still reads as bold/17pt/Calibri.
$storeName = $title.Characters(14, 14)
$storeName.Font.Name = "Calibri"
$storeName.Font.Size = 17
$storeName.Font.Bold = $true

# --- Clear the sample/demo data rows (keep header row 3 and the single
# "Food Leftover" entry in row 4; rows 5-13 become blank placeholders).
$ws.Range("B5:D5").ClearContents()
$ws.Range("A5:A13").ClearContents()

# --- Reset the view: select the title row instead of leaving the old
# scroll position / selection behind.
[void]$ws.Range("A2:F2").Select()
